$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Protokoll")

$ws.Range("C3").Value = "Fortsetzung: Recherche Responsive WPF; GUI Mockups für C#-Anwendung"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 0.2

$ws.Range("C5").Select()
